$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Every textual/semantic change in the diff lives inside a single
# paragraph: the "Live demo of the prototype" bullet. Everywhere else the
# diff is just proofErr (spell/grammar-check) markup churn with runs
# re-merged around it - no visible text changes there, nothing to do.
# -----------------------------------------------------------------------

# 1) Drop " (Jasmeen & Jose)" after "Live demo of the prototype", leaving a
#    single space before "(5 - 10 mins) ".
$rng = $d.Content
$rng.Find.Execute(" (Jasmeen & Jose) (5", $true, $false, $false, $false, $false, $true, 1, $false, " (5", 2) | Out-Null

# 2) Insert a new bold "Slide: " + non-bold "Big picture structure (Leisy)"
#    line right after "...mins) " (and the line break that already follows
#    it), landing just before the existing "Slide(s): Screenshots, verbal
#    walkthrough." line.
$anchor = $d.Content
$anchor.Find.Execute("Live demo of the prototype") | Out-Null

$scope = $d.Range($anchor.End, $d.Content.End)
$scope.Find.Execute("mins) ") | Out-Null
$scope.Collapse(0)
$scope.MoveEnd(1, 1) | Out-Null   # step over the existing <w:br/>
$scope.Collapse(0)

$scope.InsertAfter("Slide: ")
$scope.Collapse(0)

$newSlide = $d.Range($scope.End, $scope.End)
$newSlide.InsertAfter("Big picture structure (Leisy)")
$newSlide.Bold = 0
$newSlide.Collapse(0)

$brRng = $d.Range($newSlide.End, $newSlide.End)
$brRng.InsertAfter("`v")
$brRng.Bold = 1

# 3) Append "(Jasmeen and Jose)" right after "...verbal walkthrough. ".
$rng2 = $d.Content
$rng2.Find.Execute("verbal walkthrough. ") | Out-Null
$rng2.Collapse(0)
$rng2.InsertAfter("(Jasmeen and Jose)")

# 4) Word re-anchors the hidden "_GoBack" bookmark (last-edit marker) to
#    wherever editing happened; move it from the title line to right after
#    "Live demo of the prototype", matching where this edit actually landed.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$goBackAnchor = $d.Content
$goBackAnchor.Find.Execute("Live demo of the prototype") | Out-Null
$goBackTarget = $d.Range($goBackAnchor.End, $goBackAnchor.End)
$d.Bookmarks.Add("_GoBack", $goBackTarget) | Out-Null

Write-Host "done"
